$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 888.4179
$ws.Range("J17").Value = 867.6667
$ws.Range("L17").Value = 2603.0001
$ws.Range("N17").Value = -2939.0001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 64.375
$ws.Range("I55").Value = 82
$ws.Range("J55").Value = 53.8
$ws.Range("K55").Value = 82
$ws.Range("L55").Value = 53.8
$ws.Range("M55").Value = 132
$ws.Range("N55").Value = -481.8

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 10112.125
$ws.Range("I86").Value = 11128.286
$ws.Range("J86").Value = 2999
$ws.Range("K86").Value = 11128.286
$ws.Range("L86").Value = 2999
$ws.Range("M86").Value = -10005.286
$ws.Range("N86").Value = -5245

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 10112.125
$ws.Range("I89").Value = 11128.286
$ws.Range("J89").Value = 2999
$ws.Range("K89").Value = 55641.43
$ws.Range("L89").Value = 14995
$ws.Range("M89").Value = -50025.43
$ws.Range("N89").Value = -26227

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H99").Value = 180
$ws.Range("I99").Value = 180
$ws.Range("K99").Value = 540
$ws.Range("M99").Value = 958

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H118").Value = 2729.6
$ws.Range("I118").Value = 1724.5
$ws.Range("K118").Value = 5173.5
$ws.Range("M118").Value = -3516.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 5290.517
$ws.Range("I138").Value = 5093.6313
$ws.Range("J138").Value = 5664.6
$ws.Range("K138").Value = 15280.8939
$ws.Range("L138").Value = 16993.8
$ws.Range("M138").Value = -10140.8939
$ws.Range("N138").Value = -27273.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1787.7778
$ws.Range("I107").Value = 1298
$ws.Range("J107").Value = 2400
$ws.Range("K107").Value = 1298
$ws.Range("L107").Value = 2400
$ws.Range("M107").Value = 622
$ws.Range("N107").Value = -6240

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1583.7838
$ws.Range("I31").Value = 1481.0625
$ws.Range("J31").Value = 2241.2
$ws.Range("K31").Value = 1481.0625
$ws.Range("L31").Value = 2241.2
$ws.Range("M31").Value = -1186.0625
$ws.Range("N31").Value = -2831.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1583.7838
$ws.Range("I34").Value = 1481.0625
$ws.Range("J34").Value = 2241.2
$ws.Range("K34").Value = 1481.0625
$ws.Range("L34").Value = 2241.2
$ws.Range("M34").Value = -1279.0625
$ws.Range("N34").Value = -2645.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 10466.125
$ws.Range("I99").Value = 15385.8
$ws.Range("J99").Value = 2266.6667
$ws.Range("K99").Value = 15385.8
$ws.Range("L99").Value = 2266.6667
$ws.Range("M99").Value = -13887.8
$ws.Range("N99").Value = -5262.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 725.375
$ws.Range("I107").Value = 534
$ws.Range("K107").Value = 534
$ws.Range("M107").Value = 1386

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 10466.125
$ws.Range("I126").Value = 15385.8
$ws.Range("J126").Value = 2266.6667
$ws.Range("K126").Value = 46157.39999999999
$ws.Range("L126").Value = 6800.000100000001
$ws.Range("M126").Value = -43687.39999999999
$ws.Range("N126").Value = -11740.0001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2215.5908
$ws.Range("I134").Value = 2284.8948
$ws.Range("J134").Value = 1776.6666
$ws.Range("K134").Value = 6854.6844
$ws.Range("L134").Value = 5329.9998
$ws.Range("M134").Value = -4319.6844
$ws.Range("N134").Value = -10399.9998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 2206.8333
$ws.Range("I6").Value = 562.75
$ws.Range("J6").Value = 5495
$ws.Range("K6").Value = 1688.25
$ws.Range("L6").Value = 16485
$ws.Range("M6").Value = -1575.25
$ws.Range("N6").Value = -16711

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 2786.077
$ws.Range("J34").Value = 2826.2727
$ws.Range("L34").Value = 8478.8181
$ws.Range("N34").Value = -8646.8181

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 1287.25
$ws.Range("I51").Value = 383
$ws.Range("K51").Value = 1149
$ws.Range("M51").Value = -689

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 24071.428
$ws.Range("I87").Value = 11750
$ws.Range("K87").Value = 35250
$ws.Range("M87").Value = -34002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H90").Value = 24071.428
$ws.Range("I90").Value = 11750
$ws.Range("K90").Value = 105750
$ws.Range("M90").Value = -99510

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H110").Value = 16143.777
$ws.Range("I110").Value = 8431.333000000001
$ws.Range("K110").Value = 25293.999
$ws.Range("M110").Value = -21203.999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H116").Value = 30980
$ws.Range("I116").Value = 9950
$ws.Range("K116").Value = 29850
$ws.Range("M116").Value = -26408

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 15767.833
$ws.Range("I121").Value = 217.71428
$ws.Range("J121").Value = 25663.363
$ws.Range("K121").Value = 653.14284
$ws.Range("L121").Value = 76990.08900000001
$ws.Range("M121").Value = 656.85716
$ws.Range("N121").Value = -79610.08900000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 12254.125
$ws.Range("J137").Value = 12254.125
$ws.Range("L137").Value = 36762.375
$ws.Range("N137").Value = -46962.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H75").Value = 54800
$ws.Range("J75").Value = 54800
$ws.Range("L75").Value = 54800
$ws.Range("N75").Value = -56548

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H78").Value = 54800
$ws.Range("J78").Value = 54800
$ws.Range("L78").Value = 164400
$ws.Range("N78").Value = -173136

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2469.8823
$ws.Range("I122").Value = 1824
$ws.Range("K122").Value = 5472
$ws.Range("M122").Value = -3022

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H130").Value = 48650
$ws.Range("J130").Value = 48650
$ws.Range("L130").Value = 48650
$ws.Range("N130").Value = -58690

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H141").Value = 77423.71000000001
$ws.Range("J141").Value = 77423.71000000001
$ws.Range("L141").Value = 77423.71000000001
$ws.Range("N141").Value = -87783.71000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("M29").ClearContents()
$ws.Range("N29").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 7723.75
$ws.Range("I61").Value = 8383.781999999999
$ws.Range("K61").Value = 8383.781999999999
$ws.Range("M61").Value = -8181.781999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 7723.75
$ws.Range("I113").Value = 8383.781999999999
$ws.Range("K113").Value = 8383.781999999999
$ws.Range("M113").Value = -6213.781999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 57500
$ws.Range("I75").Value = 57500
$ws.Range("K75").Value = 57500
$ws.Range("M75").Value = -56564

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H78").Value = 57500
$ws.Range("I78").Value = 57500
$ws.Range("K78").Value = 172500
$ws.Range("M78").Value = -167820

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 2126.5
$ws.Range("I107").Value = 1016.1923
$ws.Range("J107").Value = 4532.1665
$ws.Range("K107").Value = 3048.5769
$ws.Range("L107").Value = 13596.4995
$ws.Range("M107").Value = -1128.5769
$ws.Range("N107").Value = -17436.4995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2348.25
$ws.Range("I132").Value = 1530.1904
$ws.Range("K132").Value = 4590.5712
$ws.Range("M132").Value = -2060.5712
